$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Lines sheet: append a new "NetworkCost" parameter row (row 9)
# ---------------------------------------------------------------------------
$wsLines = $wb.Worksheets.Item("Lines")
$wsLines.Activate()

$wsLines.Range("A9").Value = "NetworkCost"
$wsLines.Range("B9").Value = "Marginal cost of the line"
$wsLines.Range("C9").Value = "Parameter"
$wsLines.Range("D9").Value = "float"
$wsLines.Range("E9").Value = "[L]"

# Force the boolean-looking text "True" to be stored as literal text (shared
# string), matching the rest of the sheet, instead of being auto-coerced to
# a native Excel boolean value.
$fLinesBool = $wsLines.Cells.Item(9, 6)
$fLinesBool.Formula = '="True"'
$fLinesBool.Copy()
$fLinesBool.PasteSpecial(-4163)

$wsLines.Rows(9).EntireRow.Select()

# ---------------------------------------------------------------------------
# Links sheet: insert the same "NetworkCost" parameter as row 8 (pushing the
# existing HyperArcID row down to row 9)
# ---------------------------------------------------------------------------
$wsLinks = $wb.Worksheets.Item("Links")
$wsLinks.Activate()

$wsLinks.Rows(8).Insert()

# Pick up the formatting already used for manually-appended rows elsewhere in
# the workbook (Lines!A8 uses the un-bordered, non-header style) and stamp it
# across the whole new row before writing values into it.
$wsLines.Range("A8").Copy()
$wsLinks.Range("A8:F8").PasteSpecial(-4122)

$wsLinks.Range("A8").Value = "NetworkCost"
$wsLinks.Range("B8").Value = "Marginal cost of the line"
$wsLinks.Range("C8").Value = "Parameter"
$wsLinks.Range("D8").Value = "float"
$wsLinks.Range("E8").Value = "[Li]"

$fLinksBool = $wsLinks.Cells.Item(8, 6)
$fLinksBool.Formula = '="True"'
$fLinksBool.Copy()
$fLinksBool.PasteSpecial(-4163)

$wsLinks.Rows(8).EntireRow.Select()
